# Insert a new weekly record at row 178. This shifts the existing
# rows 178-292 down to 179-293 (Excel preserves formatting/styles of the
# shifted rows, including the date number-format on column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("178").Insert()

# Columns A,B,C,E,F,G,H,I,J,Q,R,T are constant for every record in this
# block, so copy them from the row directly below (which now holds what
# used to be row 178) into the freshly inserted row 178.
$ws.Range("A178").Value2 = $ws.Range("A179").Value2
$ws.Range("B178").Value2 = $ws.Range("B179").Value2
$ws.Range("C178").Value2 = $ws.Range("C179").Value2
$ws.Range("E178").Value2 = $ws.Range("E179").Value2
$ws.Range("F178").Value2 = $ws.Range("F179").Value2
$ws.Range("G178").Value2 = $ws.Range("G179").Value2
$ws.Range("H178").Value2 = $ws.Range("H179").Value2
$ws.Range("I178").Value2 = $ws.Range("I179").Value2
$ws.Range("J178").Value2 = $ws.Range("J179").Value2
$ws.Range("Q178").Value2 = $ws.Range("Q179").Value2
$ws.Range("R178").Value2 = $ws.Range("R179").Value2
$ws.Range("T178").Value2 = $ws.Range("T179").Value2

# New weekly data values for row 178.
$ws.Range("D178").Value2 = 44529
$ws.Range("K178").Value2 = "Sin especificar"
$ws.Range("L178").Value2 = "Primera Pint" + [char]0xF3 + "n"
$ws.Range("M178").Value2 = 800
$ws.Range("N178").Value2 = 24000
$ws.Range("O178").Value2 = 25000
$ws.Range("P178").Value2 = 24500
$ws.Range("S178").Value2 = 1225
